$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 30: equip ship from item condition
$ws.Range("A30").Value = "itemIsShipHeader"
$ws.Range("B30").Value = "是否是船首像"
$ws.Range("C30").Value = "item"
$ws.Range("D30").Value = "type"
$ws.Range("E30").Value = "'="
$ws.Range("E30").Style = "Normal"
$ws.Range("F30").Value = "number"
$ws.Range("G30").Value = 6

$ws.Range("A30").Select()
